# Add more blank slides to the default presentation.
#
# The deck ends with a "This slide left blank for whiteboard" slide. The
# author appended 14 additional copies of that same placeholder slide so
# the deck grows from 6 slides to 20 slides (slide 1 stays the intro/
# feature-overview slide; every slide after it is the blank whiteboard
# placeholder).
#
# We reproduce that by duplicating the current last slide (the blank
# placeholder) repeatedly, each duplicate landing right after it, until
# the deck has 20 slides total.

$p = $ppt.ActivePresentation

$targetSlideCount = 20

while ($p.Slides.Count -lt $targetSlideCount) {
    $lastIndex = $p.Slides.Count
    $p.Slides.Item($lastIndex).Duplicate() | Out-Null
}
